$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the same style used by row 5 (style index 4) to the new row 8 cells
$ws.Range("A5").Copy()
$ws.Range("A8:F8").PasteSpecial(-4122)

$values = @(
  " גם התפיסה הניהיליסטית של ריק, שקשורה באופן הדוק לאופיו – הפילוסופיה האישית שלו קיבלה השראה מאוד עמוקה מספרותה של “נרדוניה ווליה”, לדוגמה",
  "המעריצים מבינים את הדברים הללו",
  " להם את היכולת האינטלקטואלית להעריך באמת את עומק הבדיחות, להבין שהן לא סתם מצחיקות – הן אומרות משהו עמוק על החיים",
  "כתוצאה מכך, אנשים שלא אוהבים את ריק ומורטי הם באמת אידיוטים",
  "ברור שהם לא יעריכו, לדוגמה, את ההומור בביטוי הקליט האקזיסטנציאלי של ריק “וואבה לאבה דאב דאב”, שהוא בעצם רפרנס נסתר ל”אבות ובנים” של טורגנייב",
  "אני מחייך כעת רק מלדמיין את הנבערים המשתאים מגרדים בראשם בבלבול אל מול גאונותו של דן הרמון הנפרשת מולם על מסכי הטלוויזיה. אילו כסילים… כמה אני מרחם עליהם. 😂 וכן, דרך אגב, אכן יש לי קעקוע של ריק ומורטי. ולא, אתם לא יכולים לראות אותו. הוא לעיני הבחורות בלבד – וגם הן צריכות להוכיח שהן בטווח 5 נקודות אייקיו משלי (עדיף כלפי מטה) לפני כן"
)

$cols = @("A", "B", "C", "D", "E", "F")

for ($i = 0; $i -lt 6; $i++) {
  $ws.Range($cols[$i] + "8").Value = $values[$i]
}
